$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "Civil Engineering"
$ws.Range("F2").Value = "Unity University College"
$ws.Range("E3").Value = "Civil Engineering"
$ws.Range("E4").Value = "Civil Engineering"
